$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 3.503626
$ws.Range("N2").Value = 7.007251999999999
$ws.Range("O2").Value = 0.3169831060911132
$ws.Range("P2").Value = 0.2497807590538258
$ws.Range("Q2").Value = 0.05447320917266666
$ws.Range("R2").Value = 0.326839255036
$ws.Range("S2").Value = 0.3169831060911132
$ws.Range("T2").Value = 0.2497807590538258

# Row 3
$ws.Range("O3").Value = 0.122955347023628
$ws.Range("P3").Value = 0.1453320980793559
$ws.Range("S3").Value = 0.122955347023628
$ws.Range("T3").Value = 0.1453320980793559

# Row 4
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.314223
$ws.Range("N4").Value = 0.942669
$ws.Range("O4").Value = 0.02842865720977863
$ws.Range("P4").Value = 0.03360241337924067
$ws.Range("Q4").Value = 0.004885434463
$ws.Range("R4").Value = 0.043968910167
$ws.Range("S4").Value = 0.02842865720977863
$ws.Range("T4").Value = 0.03360241337924067

# Row 5
$ws.Range("M5").Value = 1.601875
$ws.Range("N5").Value = 3.20375
$ws.Range("O5").Value = 0.1449262315868481
$ws.Range("P5").Value = 0.1142009887497545
$ws.Range("Q5").Value = 0.02490541854166667
$ws.Range("R5").Value = 0.14943251125
$ws.Range("S5").Value = 0.1449262315868481
$ws.Range("T5").Value = 0.1142009887497545

# Row 6
$ws.Range("M6").Value = 3.705586
$ws.Range("N6").Value = 11.116758
$ws.Range("O6").Value = 0.335255007288947
$ws.Range("P6").Value = 0.3962683590454134
$ws.Range("Q6").Value = 0.05761321593266667
$ws.Range("R6").Value = 0.518518943394
$ws.Range("S6").Value = 0.335255007288947
$ws.Range("T6").Value = 0.3962683590454134

# Row 7
$ws.Range("M7").Value = 0.5686969999999999
$ws.Range("N7").Value = 1.706091
$ws.Range("O7").Value = 0.05145165079968518
$ws.Range("P7").Value = 0.06081538169240962
$ws.Range("Q7").Value = 0.008841911390333331
$ws.Range("R7").Value = 0.07957720251299999
$ws.Range("S7").Value = 0.05145165079968518
$ws.Range("T7").Value = 0.06081538169240962
